$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.778.40'
$ws.Range('E2').Value = '  +2.58%  '
$ws.Range('D3').Value = '2.584.38'
$ws.Range('E3').Value = '  +0.86%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '519.51'
$ws.Range('E5').Value = '  +0.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.29'
$ws.Range('E6').Value = '  -2.77%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.562'
$ws.Range('E8').Value = '  +0.15%  '
$ws.Range('D9').Value = '2.589.84'
$ws.Range('E9').Value = '  +0.55%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.53'
$ws.Range('E10').Value = '  -1.34%  '
$ws.Range('E11').Value = '  +0.08%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.329'
$ws.Range('E12').Value = '  +1.49%  '
$ws.Range('E13').Value = '  +2.50%  '
$ws.Range('D14').Value = '3.039.53'
$ws.Range('E14').Value = '  +0.79%  '
$ws.Range('D15').Value = '58.688.30'
$ws.Range('E15').Value = '  +2.42%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.30'
$ws.Range('E16').Value = '  +0.70%  '
$ws.Range('D17').Value = '2.572.58'
$ws.Range('E17').Value = '  +0.17%  '
$ws.Range('E18').Value = '  -0.80%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '336.78'
$ws.Range('E19').Value = '  +0.34%  '
$ws.Range('E20').Value = '  +0.47%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.07'
$ws.Range('E21').Value = '  -1.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.49'
$ws.Range('E22').Value = '  +3.89%  '
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('E24').Value = '  +1.23%  '
$ws.Range('E25').Value = '  +0.88%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.403'
$ws.Range('E26').Value = '  +0.54%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.998'
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.98'
$ws.Range('E28').Value = '  +0.47%  '
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('E30').Value = '  -4.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.91'
$ws.Range('E31').Value = '  -8.37%  '
$ws.Range('E32').Value = '  -0.72%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.72'
$ws.Range('E33').Value = '  +0.74%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '148.56'
$ws.Range('E34').Value = '  -0.44%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.95'
$ws.Range('E35').Value = '  -1.56%  '
$ws.Range('E36').Value = '  -1.50%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '36.63'
$ws.Range('E37').Value = '  +1.69%  '
$ws.Range('E38').Value = '  +0.91%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.820'
$ws.Range('E39').Value = '  -1.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.810'
$ws.Range('E40').Value = '  -5.62%  '
$ws.Range('E41').Value = '  -0.71%  '
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '272.39'
$ws.Range('E43').Value = '  +1.46%  '
$ws.Range('E44').Value = '  +0.85%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.587'
$ws.Range('E45').Value = '  +0.28%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0946'
$ws.Range('E46').Value = '  -0.91%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0516'
$ws.Range('E47').Value = '  -0.90%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '18.41'
$ws.Range('E48').Value = '  -2.59%  '
$ws.Range('D49').Value = '1.968.18'
$ws.Range('E49').Value = '  -0.11%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.50'
$ws.Range('E50').Value = '  -0.70%  '
$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0218'
$ws.Range('E51').Value = '  -0.74%  '
